# Add a new "Slovakia" sheet as a copy of the "Portugal" template sheet,
# placed after it, then fill in the market-specific values and make it
# the active/selected sheet (matching tabSelected moving to the new sheet).

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Copy Portugal and place the new sheet right after it (Before=$null, After=Portugal).
$portugal.Copy($null, $portugal)

# The freshly copied sheet is now the last worksheet in the workbook.
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Update the ticket reference first, then the market name, so the shared
# string table is populated in the same order as the source edit.
$slovakia.Range("B4").Value = "NGC-2930/T3223"
$slovakia.Range("B2").Value = "Slovakia Market"

# Make the new sheet the active one, selecting the Description cell.
$slovakia.Activate()
$slovakia.Range("B2").Select()
